$d = $word.ActiveDocument

# The original run read "...Farmacia GI S.A. de C.V en Zimatlán de Álvarez, Oaxaca".
# It must become "...Farmacias GI S.A. de C.V en Zimatlán de Álvarez, Oaxaca"
# (i.e. "Farmacia" -> "Farmacias"), keeping the existing bold Arial formatting.
$rng = $d.Content
$found = $rng.Find.Execute(
    " para las sucursales de Farmacia GI S.A. de C.V en Zimatlán de Álvarez, Oaxaca",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    " para las sucursales de Farmacias GI S.A. de C.V en Zimatlán de Álvarez, Oaxaca",
    2)
